$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Power PMOS" ref-numbers cell: drop Q4 from the group (Q1, Q4, Q5 -> Q1,  Q5)
$ws.Range("B21").Value = "Q1,  Q5"

# Add a new BOM row for the NMOS power transistor (Q4) that was split out of row 21
$ws.Range("B39").Value = "Q4"
$ws.Range("C39").Value = "NMOS power"
$ws.Range("D39").Value = "785-1460-1-ND"
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 0.6

# Add the Digikey hyperlink for the new part, then pick up the same visual style
# used by the other hyperlink cells in column G (reuses the existing "Hyperlink" style)
$ws.Hyperlinks.Add($ws.Range("G39"), "https://www.digikey.ca/product-detail/en/alpha-omega-semiconductor-inc/AO3434A/785-1460-1-ND/3603468")
$ws.Range("G38").Copy()
$ws.Range("G39").PasteSpecial(-4122)

# Move the selection to the newly added hyperlink cell, matching the saved view state
$ws.Range("G39").Select()
